$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort rows 2-19 alphabetically by column A (place name), keeping header row 1 fixed
$rng = $ws.Range("A1:F19")
$key1 = $ws.Range("A1")
$rng.Sort($key1, 1, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 2, [System.Reflection.Missing]::Value, 2, 1)

# Populate "Enemies during the day" (E) and "Enemies during the night" (F) with expanded enemy rosters
$ws.Range("E2").Value = "Bandit,Troll,Cannibal"
$ws.Range("F2").Value = "Wild boar,Cannibal,Fungoid,Bandit,Mad Knight,Shadow Demon,Skeleton,Specter,Vampire,Zombie"
$ws.Range("E3").Value = "Fox,Wolf,Wild boar,Cannibal,Carnivorous Plant,Centaur,Ent,Fungoid,Kobold,Toxic Spore"
$ws.Range("F3").Value = "Fox,Wolf,Wild boar,Werewolf,Fairy,Acromantula,Cannibal,Carnivorous Plant,Centaur,Ent,Fungoid,Hag,Kobold,Shadow Demon,Specter,Toxic Spore"
$ws.Range("E4").Value = "Bandit,Wolf,Cannibal,Infernal Fiend,Manticore"
$ws.Range("F4").Value = "Bandit,Wolf,Cannibal,Infernal Fiend,Mad Knight,Shadow Demon"
$ws.Range("E5").Value = "Fox,Wolf,Cannibal,Carnivorous Plant,Kobold,Toxic Spore,Wyvern"
$ws.Range("F5").Value = "Fox,Wolf,Wild boar,Fairy,Cannibal,Carnivorous Plant,Centaur,Ent,Fungoid,Kobold,Mad Knight,Shadow Demon,Succubus,Toxic Spore,Vampire"
$ws.Range("E6").Value = "Bandit,Wolf,Hag,Kobold,Toxic Spore,Zombie"
$ws.Range("F6").Value = "Bandit,Wolf,Acromantula,Fungoid,Hag,Shadow Demon,Skeleton,Specter,Toxic Spore,Zombie"
$ws.Range("E7").Value = "Wolf,Wild boar,Stone golem,Giant,Dragon,Gargoyle,Manticore,Wyvern"
$ws.Range("F7").Value = "Wolf,Wild boar,Stone golem,Giant,Orc,Goblin,Dragon,Gargoyle,Shadow Demon,Skeleton"
$ws.Range("E8").Value = "Cave troll,Zombie,Orc,Goblin,Acromantula,Basilisk,Cannibal,Dragon,Gargoyle,Skeleton,Vampire"
$ws.Range("F8").Value = "Cave troll,Werewolf,Ghost,Zombie,Orc,Goblin,Acromantula,Basilisk,Cannibal,Dragon,Gargoyle,Skeleton,Specter,Wyvern"
$ws.Range("E9").Value = "Cave troll,Dragon,Stone golem,Acromantula,Basilisk,Cannibal,Dragon,Manticore,Skeleton,Vampire"
$ws.Range("F9").Value = "Cave troll,Werewolf,Dragon,Stone golem,Acromantula,Basilisk,Cannibal,Dragon,Manticore,Skeleton,Specter,Wyvern"
$ws.Range("E10").Value = "Bandit,Wolf,Cannibal"
$ws.Range("F10").Value = "Bandit,Wolf,Cannibal,Mad Knight,Naga,Shadow Demon,Specter"
$ws.Range("E11").Value = "Fox,Fairy,Griffon,Manticore,Toxic Spore,Wyvern"
$ws.Range("F11").Value = "Fox,Fairy,Wolf,Wild boar,Griffon,Kobold,Shadow Demon,Specter,Toxic Spore,Vampire"
$ws.Range("E12").Value = "Bandit,Wolf,Naga,Wyvern"
$ws.Range("F12").Value = "Bandit,Wolf,Dragon,Hag,Naga,Shadow Demon,Specter"
$ws.Range("E13").Value = "Bandit,Wolf,Basilisk,Manticore"
$ws.Range("F13").Value = "Bandit,Wolf,Basilisk,Gargoyle,Manticore,Shadow Demon,Specter,Succubus,Vampire"
$ws.Range("E14").Value = "Mermaid,Troll,Cannibal"
$ws.Range("F14").Value = "Mermaid,Fairy,Troll,Cannibal,Mad Knight,Shadow Demon"
$ws.Range("E15").Value = "Mermaid,Griffon,Sea Serpent,The Krakken,Toxic Spore"
$ws.Range("F15").Value = "Mermaid,Fairy,Fungoid,Ent,Gargoyle,Hag,Mad Knight,Naga,Sea Serpent,Shadow Demon,The Krakken,Toxic Spore"
$ws.Range("E16").Value = "Bandit,Wolf,Griffon,Mermaid"
$ws.Range("F16").Value = "Bandit,Wolf,Fungoid,Ent,Hag,Shadow Demon,Mermaid"
$ws.Range("E17").Value = "Bandit,Wolf,Naga,Sea Serpent,The Krakken"
$ws.Range("F17").Value = "Bandit,Wolf,Mad Knight,Sea Serpent,Shadow Demon,The Krakken"
$ws.Range("E18").Value = "Fox,Mermaid,Griffon,Toxic Spore"
$ws.Range("F18").Value = "Fox,Fairy,Mermaid,Fungoid,Kobold,Manticore,Shadow Demon,Toxic Spore"
$ws.Range("E19").Value = "Bandit,Wolf,Carnivorous Plant,Fungoid,Succubus,Toxic Spore,Wyvern"
$ws.Range("F19").Value = "Bandit,Wolf,Carnivorous Plant,Fungoid,Shadow Demon,Specter,Succubus,Toxic Spore,Vampire"

# Column width adjustments (bestFit no longer applied; columns C/D get fixed widths)
$ws.Columns.Item(3).ColumnWidth = 22.35
$ws.Columns.Item(4).ColumnWidth = 14.65

# Restore the selected cell position
$ws.Range("E25").Select()
